$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh: map of cell -> new literal text.
# Values are staged through a scratch cell using a quoted-string formula
# ("=""..."""), then copied in as values-only. This keeps numeric-looking
# strings (e.g. "0.9996", "27.368.20") stored as TEXT, matching the
# existing cell contents, instead of letting Excel auto-convert them to
# numbers the way a plain Range.Value assignment would.
$updates = @(
    @{ Cell = 'D2'; Text = '27.368.20' }
    @{ Cell = 'E2'; Text = '  +0.94%  ' }
    @{ Cell = 'D3'; Text = '1.825.32' }
    @{ Cell = 'E3'; Text = '  -0.09%  ' }
    @{ Cell = 'D4'; Text = '0.9996' }
    @{ Cell = 'D5'; Text = '314.17' }
    @{ Cell = 'E5'; Text = '  +0.53%  ' }
    @{ Cell = 'D6'; Text = '0.9997' }
    @{ Cell = 'E6'; Text = '  -0.09%  ' }
    @{ Cell = 'D7'; Text = '0.4472' }
    @{ Cell = 'E7'; Text = '  -2.14%  ' }
    @{ Cell = 'D8'; Text = '0.3753' }
    @{ Cell = 'E8'; Text = '  +0.38%  ' }
    @{ Cell = 'D9'; Text = '0.07483' }
    @{ Cell = 'E9'; Text = '  +2.02%  ' }
    @{ Cell = 'D10'; Text = '0.8861' }
    @{ Cell = 'E10'; Text = '  +2.91%  ' }
    @{ Cell = 'D11'; Text = '21.01' }
    @{ Cell = 'E11'; Text = '  +0.05%  ' }
    @{ Cell = 'D12'; Text = '1.831.72' }
    @{ Cell = 'E12'; Text = '  +0.12%  ' }
    @{ Cell = 'D13'; Text = '6.757' }
    @{ Cell = 'E13'; Text = '  +0.93%  ' }
    @{ Cell = 'D14'; Text = '5.427' }
    @{ Cell = 'E14'; Text = '  +1.49%  ' }
    @{ Cell = 'D15'; Text = '93.82' }
    @{ Cell = 'E15'; Text = '  +0.91%  ' }
    @{ Cell = 'D16'; Text = '0.07119' }
    @{ Cell = 'E16'; Text = '  +0.67%  ' }
    @{ Cell = 'D17'; Text = '1.000' }
    @{ Cell = 'E17'; Text = '  -0.16%  ' }
    @{ Cell = 'D18'; Text = '0.000008786' }
    @{ Cell = 'E18'; Text = '  -0.60%  ' }
    @{ Cell = 'E19'; Text = '  -0.07%  ' }
    @{ Cell = 'D20'; Text = '15.16' }
    @{ Cell = 'E20'; Text = '  +0.94%  ' }
    @{ Cell = 'D21'; Text = '27.350.32' }
    @{ Cell = 'E21'; Text = '  +0.48%  ' }
    @{ Cell = 'D22'; Text = '5.419' }
    @{ Cell = 'E22'; Text = '  +4.34%  ' }
    @{ Cell = 'D23'; Text = '10.94' }
    @{ Cell = 'E23'; Text = '  -0.74%  ' }
    @{ Cell = 'D24'; Text = '2.057.50' }
    @{ Cell = 'E24'; Text = '  +0.13%  ' }
    @{ Cell = 'D25'; Text = '1.965' }
    @{ Cell = 'E25'; Text = '  -1.93%  ' }
    @{ Cell = 'D26'; Text = '151.21' }
    @{ Cell = 'E26'; Text = '  -0.55%  ' }
    @{ Cell = 'D27'; Text = '2.306' }
    @{ Cell = 'E27'; Text = '  +3.14%  ' }
    @{ Cell = 'E28'; Text = '  +0.01%  ' }
    @{ Cell = 'D29'; Text = '5.385' }
    @{ Cell = 'E29'; Text = '  +2.25%  ' }
    @{ Cell = 'D30'; Text = '117.92' }
    @{ Cell = 'E30'; Text = '  +0.47%  ' }
    @{ Cell = 'D31'; Text = '0.08887' }
    @{ Cell = 'E31'; Text = '  +0.05%  ' }
    @{ Cell = 'D32'; Text = '0.7845' }
    @{ Cell = 'E32'; Text = '  +3.37%  ' }
    @{ Cell = 'E33'; Text = '  +0.77%  ' }
    @{ Cell = 'D34'; Text = '4.623' }
    @{ Cell = 'E34'; Text = '  +3.32%  ' }
    @{ Cell = 'D35'; Text = '2.912' }
    @{ Cell = 'E35'; Text = '  -2.31%  ' }
    @{ Cell = 'E36'; Text = '  -0.09%  ' }
    @{ Cell = 'D37'; Text = '1.109' }
    @{ Cell = 'E37'; Text = '  +0.40%  ' }
    @{ Cell = 'E38'; Text = '  +1.11%  ' }
    @{ Cell = 'D39'; Text = '0.05301' }
    @{ Cell = 'E39'; Text = '  +0.17%  ' }
    @{ Cell = 'D40'; Text = '7.313' }
    @{ Cell = 'E40'; Text = '  +1.32%  ' }
    @{ Cell = 'D41'; Text = '0.5360' }
    @{ Cell = 'E41'; Text = '  -0.62%  ' }
    @{ Cell = 'D42'; Text = '2.857' }
    @{ Cell = 'E42'; Text = '  -1.03%  ' }
    @{ Cell = 'D43'; Text = '0.1717' }
    @{ Cell = 'E43'; Text = '  +0.38%  ' }
    @{ Cell = 'D44'; Text = '2.307' }
    @{ Cell = 'E44'; Text = '  +17.21%  ' }
    @{ Cell = 'D45'; Text = '8.669' }
    @{ Cell = 'E45'; Text = '  +0.38%  ' }
    @{ Cell = 'D46'; Text = '0.5118' }
    @{ Cell = 'E46'; Text = '  -2.21%  ' }
    @{ Cell = 'D47'; Text = '10.57' }
    @{ Cell = 'E47'; Text = '  -1.56%  ' }
    @{ Cell = 'D48'; Text = '1.699' }
    @{ Cell = 'E48'; Text = '  +1.29%  ' }
    @{ Cell = 'D49'; Text = '105.33' }
    @{ Cell = 'E49'; Text = '  -0.91%  ' }
    @{ Cell = 'D50'; Text = '0.9994' }
    @{ Cell = 'E50'; Text = '  -0.04%  ' }
    @{ Cell = 'D51'; Text = '0.06397' }
    @{ Cell = 'E51'; Text = '  +0.62%  ' }
)

$scratch = $ws.Range("Z1")

foreach ($u in $updates) {
    $scratch.Formula = "=""" + $u.Text + """"
    $scratch.Copy()
    $ws.Range($u.Cell).PasteSpecial(-4163)
}

$scratch.ClearContents()
$excel.CutCopyMode = $false